$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MOSIP_QueryLog_External")

# ---------------------------------------------------------------------------
# 1. Update the "Comments - On PD(s)" (N) and "Comments - Shrikant" O-column
#    review notes for the rows belonging to the "Registration Processor"
#    module. These replace the old `=M{row}*1.5` shared formula (which
#    evaluated to 0 because column M was blank) with literal review
#    comments/PD estimates added during the query-log review.
# ---------------------------------------------------------------------------
$ws.Range("N11").Value = "24"
$ws.Range("O11").Value = "Assumptions : RP will expose 2 apis`n1. to provide applicant info.`n2. receive updated info and incorporate it."

$ws.Range("N12").Value = "36"
$ws.Range("O12").Value = "Estimation may change after understanding overall scope of the change."

$ws.Range("N13").Value = "52"
$ws.Range("O13").Value = "Need more clarification on the requirement. Estimation may change after clarification."

$ws.Range("N14").Value = "30"

$ws.Range("N15").Value = "28"
$ws.Range("O15").Value = "Estimation may change after understanding overall scope of the change."

$ws.Range("N20").Value = "20"
$ws.Range("O20").Value = "Since the requirement is not detailed the effort may change."

$ws.Range("N40").Value = "45"

$ws.Range("N41").Value = "55"
$ws.Range("O41").Value = "Need more clarification on the requirement. Estimation may change after clarification."

$ws.Range("N42").Value = "12"
$ws.Range("O42").Value = "Change algorithm from lavenstine distance to phonetic and soundex match"

# ---------------------------------------------------------------------------
# 2. Filter the log down to the "Registration Processor" module rows (column
#    D) across the full data range (the sheet grew from 34 to 53 data rows
#    since the filter was first applied), hiding every non-matching row.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A2:H53").AutoFilter(4, @("Registration Processor"), 7)

# Keep the hidden "_FilterDatabase" defined name in sync with the new
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "MOSIP_QueryLog_External!_FilterDatabase") {
        $n.RefersTo = "=MOSIP_QueryLog_External!`$A`$2:`$H`$53"
    }
}

# ---------------------------------------------------------------------------
# 3. Restore the working selection to the last-reviewed cell.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("O13").Select()
